$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Search" sheet had a redundant column L ("STATUS") value of "PASS" copy-pasted
# into almost every data row (rows 2-31). Clear that whole block in a single range
# operation so the writer collapses each row's declared column span (and drops the
# now-empty trailing <c> elements) exactly like Excel does when you select L2:L31 and
# hit Delete.
$ws.Range("L2:L31").ClearContents()

# Update the view state left behind by the editor: scrolled over to show column H,
# with K6 as the active selected cell.
$ws.Range("K6").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
